$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column BF = column 58. Rows 2-31 hold the per-game "Date" value which was
# mistakenly written as "5-19-2011-12" (an artifact of how the NBA stats
# feed rendered the date). Correct it to ISO "2012-05-19", keeping the
# cell a plain text value (not an auto-converted date serial).
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    if ($cell.Value2 -eq "5-19-2011-12") {
        $cell.NumberFormat = "@"
        $cell.Value2 = "2012-05-19"
        $cell.Style = "Normal"
    }
}
